# GenX FX trading signals refresh:
#  - "Active Signals": new snapshot with only 2 live signals (rows 4-7 removed).
#  - "Summary Dashboard": refreshed aggregate counters.
#  - "Signal History": rows reshuffled with the newest market data + statuses.
#
# Helper: writes a value as literal TEXT (never let Excel's "looks like a
# number/percentage" auto-detection kick in), while leaving the cell's
# existing style/format (borders, fill, font) completely untouched -- we
# stamp the original format back on top immediately after the write.
function Set-TextValue {
    param($range, [string]$text)
    $fmtSource = $range.Worksheet.Cells.Item($range.Row, $range.Column)
    $range.Value = "'" + $text
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Active Signals"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Active Signals")

# Row 2 (XAUAUD BUY) refreshed prices / confidence.
$ws1.Range("A2").Value = "2025-07-28 21:28"
$ws1.Range("D2").Value = 4064.91481
$ws1.Range("E2").Value = 4064.91121
$ws1.Range("F2").Value = 4064.91896
$ws1.Range("G2").Value = 0.01
Set-TextValue $ws1.Range("H2") "84.0%"

# Row 3 becomes a brand-new SELL signal (XAUGBP) replacing the old XAUCAD BUY.
$ws1.Range("A3").Value = "2025-07-28 20:55"
$ws1.Range("B3").Value = "XAUGBP"
$ws1.Range("C3").Value = "SELL"
$ws1.Range("D3").Value = 2109.70362
$ws1.Range("E3").Value = 2109.7061
$ws1.Range("F3").Value = 2109.69605
$ws1.Range("G3").Value = 0.09
Set-TextValue $ws1.Range("H3") "84.0%"
$ws1.Range("I3").Value = 3.04
# J3 stays "Active" (unchanged).

# SELL signals get a light-red fill (mirrors the existing light-green BUY fill).
$ws1.Range("C3").Interior.Color = 13551615

# Rows 4-7 (XAUCAD/XAUEUR/XAUEUR/XAUUSD/XAUCHF) are gone now -- only 2 active
# signals remain.
$ws1.Rows("4:7").Delete()

# ---------------------------------------------------------------------------
# Sheet 2: "Summary Dashboard"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary Dashboard")

$ws2.Range("B4").Value = 2
$ws2.Range("B5").Value = 7
$ws2.Range("B6").Value = 8
Set-TextValue $ws2.Range("B7") "83.5%"
Set-TextValue $ws2.Range("B8") "1.93"
Set-TextValue $ws2.Range("B9") "2025-07-28 21:07:35"

# ---------------------------------------------------------------------------
# Sheet 3: "Signal History"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Signal History")

function Set-HistoryRow {
    param($ws, [int]$row, [string]$ts, [string]$sym, [string]$sig,
          [double]$entry, [double]$sl, [double]$tp, [double]$lots,
          [double]$conf, [double]$rr, [string]$status)
    $ws.Range("A$row").Value = $ts
    $ws.Range("B$row").Value = $sym
    $ws.Range("C$row").Value = $sig
    $ws.Range("D$row").Value = $entry
    $ws.Range("E$row").Value = $sl
    $ws.Range("F$row").Value = $tp
    $ws.Range("G$row").Value = $lots
    $ws.Range("H$row").Value = $conf
    $ws.Range("I$row").Value = $rr
    $ws.Range("J$row").Value = $status
}

Set-HistoryRow $ws3 2  "2025-07-28 21:19" "XAUUSD" "BUY"  2638.81797 2638.81536 2638.82472 0.03 0.91 2.59 "Filled"
Set-HistoryRow $ws3 3  "2025-07-28 20:56" "NZDUSD" "SELL" 0.5863     0.58862    0.58136    0.1  0.85 2.12 "Pending"
Set-HistoryRow $ws3 4  "2025-07-28 20:43" "EURUSD" "SELL" 1.10395    1.10659    1.09987    0.02 0.78 1.54 "Filled"
Set-HistoryRow $ws3 5  "2025-07-28 20:54" "XAUCAD" "SELL" 3602.6381  3602.64162 3602.63223 0.1  0.85 1.67 "Pending"
Set-HistoryRow $ws3 6  "2025-07-28 21:28" "XAUAUD" "BUY"  4064.91481 4064.91121 4064.91896 0.01 0.84 1.15 "Active"
Set-HistoryRow $ws3 7  "2025-07-28 20:55" "XAUGBP" "SELL" 2109.70362 2109.7061  2109.69605 0.09 0.84 3.04 "Active"
Set-HistoryRow $ws3 8  "2025-07-28 20:40" "XAUCAD" "BUY"  3637.04486 3637.0413  3637.05461 0.09 0.77 2.75 "Pending"
Set-HistoryRow $ws3 9  "2025-07-28 21:33" "XAUCHF" "BUY"  2330.19431 2330.19843 2330.18961 0.09 0.79 1.14 "Filled"
Set-HistoryRow $ws3 10 "2025-07-28 21:19" "XAUUSD" "BUY"  2654.13881 2654.13442 2654.14534 0.06 0.75 1.48 "Pending"
Set-HistoryRow $ws3 11 "2025-07-28 20:52" "XAUEUR" "BUY"  2414.83832 2414.84059 2414.82938 0.03 0.77 3.94 "Filled"
Set-HistoryRow $ws3 12 "2025-07-28 20:59" "NZDUSD" "BUY"  0.58938    0.59428    0.58413    0.02 0.89 1.07 "Pending"
Set-HistoryRow $ws3 13 "2025-07-28 21:25" "EURUSD" "BUY"  1.10743    1.10362    1.1122     0.04 0.77 1.25 "Pending"
Set-HistoryRow $ws3 14 "2025-07-28 21:30" "XAUEUR" "BUY"  2412.942   2412.93763 2412.94933 0.07000000000000001 0.9  1.68 "Filled"
Set-HistoryRow $ws3 15 "2025-07-28 21:35" "USDCAD" "BUY"  1.3615     1.36633    1.35425    0.02 0.87 1.5  "Pending"
Set-HistoryRow $ws3 16 "2025-07-28 21:07" "USDCHF" "BUY"  0.88156    0.87934    0.88597    0.03 0.9399999999999999 1.99 "Filled"
